# Update the "取得日時" (fetched at) timestamp column for the appended
# batch of rows on the "ランサーズ" sheet to reflect the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$newTimestamp = "2025-11-07 18:32:38"

for ($row = 2; $row -le 13; $row++) {
    $ws.Cells.Item($row, 1).Value = $newTimestamp
}
